$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F holds "dSF" values (row 1 is header). Update rows 2-8 per the new data pull.
$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = -4
